$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp and country-name reorderings
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 22:07"
$ws.Range("A86").Value = "Costa Rica"
$ws.Range("A87").Value = "Tayikistan"
$ws.Range("A88").Value = "Bulgaria"
$ws.Range("A89").Value = "Bosnia y Herzegovina"
$ws.Range("A90").Value = "Gabon"
$ws.Range("A184").Value = "Lesoto"
$ws.Range("A185").Value = "Seychelles"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# Update refreshed numeric data (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes)
$ws.Range("B4").Value = 3198109
$ws.Range("C4").Value = 39177
$ws.Range("D4").Value = 1408595
$ws.Range("E4").Value = 1654062
$ws.Range("G4").Value = 590
$ws.Range("H4").Value = 135452
$ws.Range("B16").Value = 238339
$ws.Range("C16").Value = 13674
$ws.Range("D16").Value = 113061
$ws.Range("E16").Value = 121558
$ws.Range("G16").Value = 118
$ws.Range("H16").Value = 3720
$ws.Range("B19").Value = 199077
$ws.Range("C19").Value = 312
$ws.Range("E19").Value = 6353
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 9124
$ws.Range("B70").Value = 11750
$ws.Range("C70").Value = 246
$ws.Range("D70").Value = 5752
$ws.Range("E70").Value = 5919
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 79
$ws.Range("B84").Value = 6973
$ws.Range("C84").Value = 199
$ws.Range("E84").Value = 4423
$ws.Range("B86").Value = 6485
$ws.Range("C86").Value = 649
$ws.Range("D86").Value = 2023
$ws.Range("E86").Value = 4437
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 25
$ws.Range("B87").Value = 6410
$ws.Range("C87").Value = 46
$ws.Range("D87").Value = 5067
$ws.Range("E87").Value = 1289
$ws.Range("H87").Value = 54
$ws.Range("B88").Value = 6342
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 3166
$ws.Range("E88").Value = 2917
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 259
$ws.Range("B89").Value = 6086
$ws.Range("C89").Value = 217
$ws.Range("D89").Value = 2815
$ws.Range("E89").Value = 3057
$ws.Range("G89").Value = 5
$ws.Range("H89").Value = 214
$ws.Range("B90").Value = 5871
$ws.Range("D90").Value = 2682
$ws.Range("E90").Value = 3143
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 46
